# Amend wording to set dates question
#
# Original:
#   "Date the transfer was first discussed with outgoing trust"
# New:
#   "Date the transfer was first discussed with incoming or outgoing trust"
# split across three runs (same bold formatting throughout):
#   1) "Date the transfer was first discussed with "
#   2) "incoming or "
#   3) "outgoing trust"

$d = $word.ActiveDocument

$oldText = "Date the transfer was first discussed with outgoing trust"
$insertText = "incoming or "
$anchorText = "Date the transfer was first discussed with "

# Locate the target sentence.
$target = $d.Content
$found = $target.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $sentenceStart = $target.Start

    $insertPos = $sentenceStart + $anchorText.Length

    # Insert the new wording in the middle of the run.
    $insertionPoint = $d.Range($insertPos, $insertPos)
    $insertionPoint.InsertBefore($insertText)

    # Force the inserted text (and the run it was inserted into) to split
    # off from its neighbours into their own runs by toggling a character
    # format away from, then back to, the shared formatting.
    $newRange = $d.Range($insertPos, $insertPos + $insertText.Length)
    $newRange.Bold = 0
    $newRange.Bold = 1
}
